$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 811.8946999999999
$ws.Range("I32").Value = 366.6
$ws.Range("J32").Value = 970.9286
$ws.Range("K32").Value = 366.6
$ws.Range("L32").Value = 970.9286
$ws.Range("M32").Value = -40.60000000000002
$ws.Range("N32").Value = -1622.9286

$ws.Range("H68").Value = 37023.75
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 37023.75
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 37023.75
$ws.Range("N68").Value = -38521.75

$ws.Range("H71").Value = 37023.75
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 37023.75
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 111071.25
$ws.Range("N71").Value = -118559.25

$ws.Range("H75").Value = 72055.5
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 72055.5
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 72055.5
$ws.Range("M75").Value = ""
$ws.Range("N75").Value = -73927.5

$ws.Range("H78").Value = 72055.5
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 72055.5
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 216166.5
$ws.Range("M78").Value = ""
$ws.Range("N78").Value = -225526.5

$ws.Range("H111").Value = 500005000
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 500005000
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 1500015000
$ws.Range("M111").Value = ""
$ws.Range("N111").Value = -1500021134

$ws.Range("H121").Value = 1656.75
$ws.Range("I121").Value = 931.6667
$ws.Range("J121").Value = 1824.0769
$ws.Range("K121").Value = 2795.0001
$ws.Range("L121").Value = 5472.2307
$ws.Range("M121").Value = -1048.0001
$ws.Range("N121").Value = -8966.2307

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 850.1163
$ws.Range("I2").Value = 692.45715
$ws.Range("J2").Value = 1539.875
$ws.Range("K2").Value = 692.45715
$ws.Range("L2").Value = 1539.875
$ws.Range("M2").Value = -579.45715
$ws.Range("N2").Value = ""

$ws.Range("H75").Value = 110000
$ws.Range("I75").Value = 20000
$ws.Range("J75").Value = 200000
$ws.Range("K75").Value = 20000
$ws.Range("L75").Value = 200000
$ws.Range("M75").Value = -19126
$ws.Range("N75").Value = -201748

$ws.Range("H78").Value = 110000
$ws.Range("I78").Value = 20000
$ws.Range("J78").Value = 200000
$ws.Range("K78").Value = 60000
$ws.Range("L78").Value = 600000
$ws.Range("M78").Value = -55632
$ws.Range("N78").Value = -608736

$ws.Range("H116").Value = 850.1163
$ws.Range("I116").Value = 692.45715
$ws.Range("J116").Value = 1539.875
$ws.Range("K116").Value = 692.45715
$ws.Range("L116").Value = 1539.875
$ws.Range("M116").Value = 1601.54285
$ws.Range("N116").Value = ""

$ws.Range("H132").Value = 1903.1034
$ws.Range("I132").Value = 1240.625
$ws.Range("J132").Value = 2718.4614
$ws.Range("K132").Value = 3721.875
$ws.Range("L132").Value = 8155.3842
$ws.Range("M132").Value = -1191.875
$ws.Range("N132").Value = -13215.3842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 850.1163
$ws.Range("I3").Value = 692.45715
$ws.Range("J3").Value = 1539.875
$ws.Range("K3").Value = 692.45715
$ws.Range("L3").Value = 1539.875
$ws.Range("M3").Value = -578.45715
$ws.Range("N3").Value = ""

$ws.Range("H20").Value = 2277.8823
$ws.Range("I20").Value = 1940.3077
$ws.Range("J20").Value = 3375
$ws.Range("K20").Value = 1940.3077
$ws.Range("L20").Value = 3375
$ws.Range("M20").Value = -1693.3077
$ws.Range("N20").Value = -3869

$ws.Range("H51").Value = 33774
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 33774
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 33774
$ws.Range("N51").Value = -34756

$ws.Range("H76").Value = 88888
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 88888
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 88888
$ws.Range("N76").Value = -89518

$ws.Range("H79").Value = 88888
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 88888
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 88888
$ws.Range("N79").Value = -91072

$ws.Range("H107").Value = 961.1
$ws.Range("I107").Value = 876.375
$ws.Range("J107").Value = 1300
$ws.Range("K107").Value = 876.375
$ws.Range("L107").Value = 1300
$ws.Range("M107").Value = 1043.625
$ws.Range("N107").Value = ""

$ws.Range("H134").Value = 1036.6333
$ws.Range("I134").Value = 879.22644
$ws.Range("J134").Value = 2228.4285
$ws.Range("K134").Value = 2637.67932
$ws.Range("L134").Value = 6685.2855
$ws.Range("M134").Value = -102.6793200000002
$ws.Range("N134").Value = -11755.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1387.25
$ws.Range("I16").Value = 1548.7778
$ws.Range("J16").Value = 902.6667
$ws.Range("K16").Value = 1548.7778
$ws.Range("L16").Value = 902.6667
$ws.Range("M16").Value = -1261.7778
$ws.Range("N16").Value = -1476.6667

$ws.Range("H47").Value = 30035.5
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 30035.5
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 30035.5
$ws.Range("N47").Value = -31167.5

$ws.Range("H58").Value = 2559.3635
$ws.Range("I58").Value = 2625.75
$ws.Range("J58").Value = 2521.4285
$ws.Range("K58").Value = 2625.75
$ws.Range("L58").Value = 2521.4285
$ws.Range("M58").Value = -2422.75
$ws.Range("N58").Value = -2927.4285

$ws.Range("H76").Value = 1380
$ws.Range("I76").Value = 1380
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 1380
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -1065

$ws.Range("H79").Value = 1380
$ws.Range("I79").Value = 1380
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 1380
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -288

$ws.Range("H86").Value = 4433.7393
$ws.Range("I86").Value = 4991.5835
$ws.Range("J86").Value = 3825.182
$ws.Range("K86").Value = 4991.5835
$ws.Range("L86").Value = 3825.182
$ws.Range("M86").Value = -3868.5835
$ws.Range("N86").Value = -6071.182

$ws.Range("H89").Value = 4433.7393
$ws.Range("I89").Value = 4991.5835
$ws.Range("J89").Value = 3825.182
$ws.Range("K89").Value = 24957.9175
$ws.Range("L89").Value = 19125.91
$ws.Range("M89").Value = -19341.9175
$ws.Range("N89").Value = -30357.91

$ws.Range("H113").Value = 1387.25
$ws.Range("I113").Value = 1548.7778
$ws.Range("J113").Value = 902.6667
$ws.Range("K113").Value = 1548.7778
$ws.Range("L113").Value = 902.6667
$ws.Range("M113").Value = 621.2221999999999
$ws.Range("N113").Value = -5242.6667

$ws.Range("H136").Value = 2559.3635
$ws.Range("I136").Value = 2625.75
$ws.Range("J136").Value = 2521.4285
$ws.Range("K136").Value = 7877.25
$ws.Range("L136").Value = 7564.2855
$ws.Range("M136").Value = -5327.25
$ws.Range("N136").Value = -12664.2855

$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 617.05884
$ws.Range("I86").Value = 423.07693
$ws.Range("J86").Value = 1247.5
$ws.Range("K86").Value = 1269.23079
$ws.Range("L86").Value = 3742.5
$ws.Range("M86").Value = -83.23079000000007
$ws.Range("N86").Value = -6114.5

$ws.Range("H89").Value = 617.05884
$ws.Range("I89").Value = 423.07693
$ws.Range("J89").Value = 1247.5
$ws.Range("K89").Value = 3807.69237
$ws.Range("L89").Value = 11227.5
$ws.Range("M89").Value = 2120.30763
$ws.Range("N89").Value = -23083.5

$ws.Range("H107").Value = 681216.5
$ws.Range("I107").Value = 780
$ws.Range("J107").Value = 908028.7
$ws.Range("K107").Value = 2340
$ws.Range("L107").Value = 2724086.1
$ws.Range("M107").Value = -420
$ws.Range("N107").Value = -2727926.1

$ws.Range("H109").Value = 5402.25
$ws.Range("I109").Value = 4574
$ws.Range("J109").Value = 5899.2
$ws.Range("K109").Value = 13722
$ws.Range("L109").Value = 17697.6
$ws.Range("M109").Value = -12682
$ws.Range("N109").Value = ""

$ws.Range("H131").Value = 827.62
$ws.Range("I131").Value = 351.66666
$ws.Range("J131").Value = 858
$ws.Range("K131").Value = 1054.99998
$ws.Range("L131").Value = 2574
$ws.Range("M131").Value = 3985.00002
$ws.Range("N131").Value = -12654

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 980.06665
$ws.Range("I113").Value = 964.3570999999999
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 964.3570999999999
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 1205.6429
$ws.Range("N113").Value = ""

$ws.Range("H135").Value = 38663.75
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 38663.75
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 38663.75
$ws.Range("N135").Value = -48803.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 16999.666
$ws.Range("I45").Value = 8999.5
$ws.Range("J45").Value = 33000
$ws.Range("K45").Value = 8999.5
$ws.Range("L45").Value = 33000
$ws.Range("M45").Value = -8592.5
$ws.Range("N45").Value = ""

$ws.Range("H61").Value = 2002.5
$ws.Range("I61").Value = 3000
$ws.Range("J61").Value = 1005
$ws.Range("K61").Value = 3000
$ws.Range("L61").Value = 1005
$ws.Range("M61").Value = -2798
$ws.Range("N61").Value = -1409

$ws.Range("H113").Value = 2002.5
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 1005
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 1005
$ws.Range("M113").Value = -830
$ws.Range("N113").Value = -5345

$ws.Range("H127").Value = 49153.89
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 49153.89
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 49153.89
$ws.Range("N127").Value = -59073.89

$ws.Range("H133").Value = 32166.25
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 32166.25
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 32166.25
$ws.Range("N133").Value = -37226.25

$ws.Range("H136").Value = 17547030
$ws.Range("I136").Value = 3109.5625
$ws.Range("J136").Value = 111114610
$ws.Range("K136").Value = 9328.6875
$ws.Range("L136").Value = 333343830
$ws.Range("M136").Value = -6778.6875
$ws.Range("N136").Value = -333348930

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 43539
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 43539
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 43539
$ws.Range("N80").Value = -45535

$ws.Range("H83").Value = 43539
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 43539
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 130617
$ws.Range("N83").Value = -140601
